$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet has header (row 1), a summary/example row (row 2), then a large
# block of completely empty rows before the real data resumes. Two of those
# blank rows are removed, shifting all the data below up by two rows
# (dimension A1:J290 -> A1:J288, data blocks 163-189/251-290 -> 161-187/249-288).
$ws.Rows("3:4").Delete() | Out-Null

# Row 2 (the sample/demo record) gets a couple of cells reformatted to match
# the border style already used by its neighboring cells (no bottom border).
$ws.Range("B2").Copy()
$ws.Range("A2").PasteSpecial(-4122)
$ws.Range("D2").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Update the sample record's content:
#  - "Titre" (column D) changes from "Etudiant" to "Visiteur"
#  - "TypeTitreDeSejour" (column B) text had an accented typo fixed
#  - "NumEtrangerVisa" (column E) gets a new sample number
$ws.Range("D2").Value = "Visiteur"
$ws.Range("B2").Value = "RenouvellementDeTitreSejour"
$ws.Range("E2").Value = 7703039887

# Keep the _FilterDatabase defined name in sync with the new data extent.
foreach ($n in $wb.Names) {
    if ($n.Name -eq "JDD!_FilterDatabase") {
        $n.RefersTo = "=JDD!`$A`$1:`$J`$288"
    }
}

# Restore the selection recorded in the saved workbook.
$ws.Range("H33").Select() | Out-Null
